$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record for Berenjena @ Terminal La Palmera de La Serena
# is inserted as row 318, pushing the existing rows 318-325 down to 319-326.
$ws.Rows.Item(318).Insert()

# Populate the newly inserted row 318 with the new weekly data.
# Columns A,B,C,E,F,G,H,I,N,O,Q,R carry the same values the series already
# used (same market/category/unit-of-sale), only the date and the
# volume/price figures are new.
$ws.Cells.Item(318, 1).Value = 8
$ws.Cells.Item(318, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(318, 3).Value = "Coquimbo"
$ws.Cells.Item(318, 4).Value = 45239
$ws.Cells.Item(318, 5).Value = 4
$ws.Cells.Item(318, 6).Value = 100112001
$ws.Cells.Item(318, 7).Value = "Berenjena"
$ws.Cells.Item(318, 8).Value = "Sin especificar"
$ws.Cells.Item(318, 9).Value = "Primera"
$ws.Cells.Item(318, 10).Value = 440
$ws.Cells.Item(318, 11).Value = 11000
$ws.Cells.Item(318, 12).Value = 12000
$ws.Cells.Item(318, 13).Value = 11500
$ws.Cells.Item(318, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(318, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(318, 16).Value = 230
$ws.Cells.Item(318, 17).Value = 50
$ws.Cells.Item(318, 18).Value = "Hortaliza"
